$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.409.79"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "1.894.94"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = "237.61"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "0.4845"
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("D8").Value = "0.2908"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").Value = "0.06615"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").Value = "1.912.74"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").Value = "16.97"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "0.07328"
$ws.Range("D13").Value = "5.174"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "87.87"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "0.6631"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").Value = "30.393.76"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").Value = "13.48"
$ws.Range("D18").Value = "0.000007783"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("D19").Value = "0.9981"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").Value = "5.434"
$ws.Range("E20").Value = "  +4.19%  "
$ws.Range("D21").Value = "2.137.81"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").Value = "0.9976"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "194.39"
$ws.Range("E23").Value = "  -4.18%  "
$ws.Range("D24").Value = "6.198"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "9.351"
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("D26").Value = "165.13"
$ws.Range("E26").Value = "  +2.88%  "
$ws.Range("D27").Value = "18.21"
$ws.Range("E27").Value = "  -3.81%  "
$ws.Range("D28").Value = "1.944"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").Value = "1.452"
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("D30").Value = "4.317"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").Value = "0.09165"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").Value = "4.054"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "0.05093"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("D34").Value = "1.154"
$ws.Range("E34").Value = "  +2.86%  "
$ws.Range("D35").Value = "0.7281"
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("D36").Value = "2.695"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").Value = "0.01793"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("D38").Value = "2.651"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "0.9214"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").Value = "2.082"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "106.28"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").Value = "0.4329"
$ws.Range("E42").Value = "  -3.58%  "
$ws.Range("D43").Value = "5.875"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").Value = "7.562"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("E46").Value = "  -4.94%  "
$ws.Range("D47").Value = "1.567"
$ws.Range("E47").Value = "  +9.02%  "
$ws.Range("D48").Value = "65.20"
$ws.Range("E48").Value = "  -10.53%  "
$ws.Range("D49").Value = "8.948"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").Value = "34.01"
$ws.Range("E50").Value = "  -5.94%  "
$ws.Range("D51").Value = "0.05763"
$ws.Range("E51").Value = "  -3.09%  "
